# The sheet lists one Comtrade export row per reporter country (date, value,
# reporter, partner). A missing row for "Ghana" needs to be inserted between
# "Germany" (row 59) and "Kiribati" (row 59, soon to become row 60), pushing
# every following row down by one. The sheet's used-range grows from
# A1:D162 to A1:D163 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 59, shifting rows 59..162 down to 60..163.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row with the Ghana export record.
$ws.Cells.Item(59, 1).Value = 44197
$ws.Cells.Item(59, 2).Value = 14823232488.086
$ws.Cells.Item(59, 3).Value = "Ghana"
$ws.Cells.Item(59, 4).Value = "World"
